# Apply updated TPM-derived values to the Col1a1-Itga2 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "G2" = 28.56536566666667
    "H2" = 85.69609700000001
    "I2" = 0.02097368575335975
    "J2" = 0.02097368575335974
    "M2" = 3.339352
    "N2" = 10.018056
    "O2" = 0.6054960700393903
    "P2" = 0.6054960700393903
    "Q2" = 95.38981096971469
    "R2" = 858.5082987274322
    "S2" = 0.01269948429790048
    "T2" = 0.01269948429790048
    "G3" = 28.56536566666667
    "H3" = 85.69609700000001
    "I3" = 0.02097368575335975
    "J3" = 0.02097368575335974
    "O3" = 0.2540955070726236
    "P3" = 0.2540955070726236
    "Q3" = 40.03018943844601
    "R3" = 360.2717049460141
    "S3" = 0.005329319316681807
    "T3" = 0.005329319316681806
    "G4" = 28.56536566666667
    "H4" = 85.69609700000001
    "I4" = 0.02097368575335975
    "J4" = 0.02097368575335974
    "K4" = 2
    "L4" = 0.6666666666666666
    "M4" = 0.1338136666666667
    "N4" = 0.401441
    "O4" = 0.02426328499787613
    "P4" = 0.02426328499787612
    "Q4" = 3.822436319530779
    "R4" = 34.40192687577701
    "S4" = 0.0005088905148896618
    "T4" = 0.0005088905148896616
    "G5" = 28.56536566666667
    "H5" = 85.69609700000001
    "I5" = 0.02097368575335975
    "J5" = 0.02097368575335974
    "M5" = 0.6405483333333334
    "N5" = 1.921645
    "O5" = 0.11614513789011
    "P5" = 0.11614513789011
    "Q5" = 18.29749736884056
    "R5" = 164.677476319565
    "S5" = 0.002435991623887804
    "T5" = 0.002435991623887804
    "I6" = 0.9460766961189575
    "J6" = 0.9460766961189573
    "M6" = 3.339352
    "N6" = 10.018056
    "O6" = 0.6054960700393903
    "P6" = 0.6054960700393903
    "Q6" = 4302.823941718643
    "R6" = 38725.41547546779
    "S6" = 0.5728457214558793
    "T6" = 0.5728457214558792
    "I7" = 0.9460766961189575
    "J7" = 0.9460766961189573
    "O7" = 0.2540955070726236
    "P7" = 0.2540955070726236
    "S7" = 0.240393837829939
    "T7" = 0.2403938378299388
    "I8" = 0.9460766961189575
    "J8" = 0.9460766961189573
    "K8" = 2
    "L8" = 0.6666666666666666
    "M8" = 0.1338136666666667
    "N8" = 0.401441
    "O8" = 0.02426328499787613
    "P8" = 0.02426328499787612
    "Q8" = 172.4216700313388
    "R8" = 1551.795030282049
    "S8" = 0.02295492850778331
    "T8" = 0.0229549285077833
    "I9" = 0.9460766961189575
    "J9" = 0.9460766961189573
    "M9" = 0.6405483333333334
    "N9" = 1.921645
    "O9" = 0.11614513789011
    "P9" = 0.11614513789011
    "Q9" = 825.3597417986007
    "R9" = 7428.237676187407
    "S9" = 0.1098822083253561
    "T9" = 0.109882208325356
    "G10" = 0.115045
    "H10" = 0.345135
    "I10" = 0.00008447004339632664
    "J10" = 0.00008447004339632662
    "M10" = 3.339352
    "N10" = 10.018056
    "O10" = 0.6054960700393903
    "P10" = 0.6054960700393903
    "Q10" = 0.3841757508400001
    "R10" = 3.457581757560001
    "S10" = 0.00005114627931253253
    "T10" = 0.00005114627931253252
    "G11" = 0.115045
    "H11" = 0.345135
    "I11" = 0.00008447004339632664
    "J11" = 0.00008447004339632662
    "O11" = 0.2540955070726236
    "P11" = 0.2540955070726236
    "Q11" = 0.16121877093
    "R11" = 1.45096893837
    "S11" = 0.00002146345850923614
    "T11" = 0.00002146345850923613
    "G12" = 0.115045
    "H12" = 0.345135
    "I12" = 0.00008447004339632664
    "J12" = 0.00008447004339632662
    "K12" = 2
    "L12" = 0.6666666666666666
    "M12" = 0.1338136666666667
    "N12" = 0.401441
    "O12" = 0.02426328499787613
    "P12" = 0.02426328499787612
    "Q12" = 0.01539459328166667
    "R12" = 0.138551339535
    "S12" = 0.000002049520736708037
    "T12" = 0.000002049520736708037
    "G13" = 0.115045
    "H13" = 0.345135
    "I13" = 0.00008447004339632664
    "J13" = 0.00008447004339632662
    "M13" = 0.6405483333333334
    "N13" = 1.921645
    "O13" = 0.11614513789011
    "P13" = 0.11614513789011
    "Q13" = 0.07369188300833335
    "R13" = 0.6632269470750001
    "S13" = 0.000009810784837849935
    "T13" = 0.000009810784837849933
    "G14" = 43.90798866666668
    "H14" = 131.723966
    "I14" = 0.03223877359397412
    "J14" = 0.0322387735939741
    "M14" = 3.339352
    "N14" = 10.018056
    "O14" = 0.6054960700393903
    "P14" = 0.6054960700393903
    "Q14" = 146.6242297700107
    "R14" = 1319.618067930096
    "S14" = 0.019520450714041
    "T14" = 0.01952045071404099
    "G15" = 43.90798866666668
    "H15" = 131.723966
    "I15" = 0.03223877359397412
    "J15" = 0.0322387735939741
    "O15" = 0.2540955070726236
    "P15" = 0.2540955070726236
    "Q15" = 61.53063554998802
    "R15" = 553.7757199498922
    "S15" = 0.008191727523760363
    "T15" = 0.008191727523760357
    "G16" = 43.90798866666668
    "H16" = 131.723966
    "I16" = 0.03223877359397412
    "J16" = 0.0322387735939741
    "K16" = 2
    "L16" = 0.6666666666666666
    "M16" = 0.1338136666666667
    "N16" = 0.401441
    "O16" = 0.02426328499787613
    "P16" = 0.02426328499787612
    "Q16" = 5.875488959445113
    "R16" = 52.87940063500601
    "S16" = 0.0007822185516925972
    "T16" = 0.0007822185516925968
    "G17" = 43.90798866666668
    "H17" = 131.723966
    "I17" = 0.03223877359397412
    "J17" = 0.0322387735939741
    "M17" = 0.6405483333333334
    "N17" = 1.921645
    "O17" = 0.11614513789011
    "P17" = 0.11614513789011
    "Q17" = 28.12518896045223
    "R17" = 253.12670064407
    "S17" = 0.003744376804480162
    "T17" = 0.00374437680448016
    "G18" = 0.693788
    "H18" = 2.081364
    "I18" = 0.0005094032984297506
    "J18" = 0.0005094032984297505
    "M18" = 3.339352
    "N18" = 10.018056
    "O18" = 0.6054960700393903
    "P18" = 0.6054960700393903
    "Q18" = 2.316802345376
    "R18" = 20.851221108384
    "S18" = 0.0003084416952643167
    "T18" = 0.0003084416952643167
    "G19" = 0.693788
    "H19" = 2.081364
    "I19" = 0.0005094032984297506
    "J19" = 0.0005094032984297505
    "O19" = 0.2540955070726236
    "P19" = 0.2540955070726236
    "Q19" = 0.9722425889520001
    "R19" = 8.750183300568001
    "S19" = 0.0001294370894189745
    "T19" = 0.0001294370894189744
    "G20" = 0.693788
    "H20" = 2.081364
    "I20" = 0.0005094032984297506
    "J20" = 0.0005094032984297505
    "K20" = 2
    "L20" = 0.6666666666666666
    "M20" = 0.1338136666666667
    "N20" = 0.401441
    "O20" = 0.02426328499787613
    "P20" = 0.02426328499787612
    "Q20" = 0.09283831616933334
    "R20" = 0.835544845524
    "S20" = 0.00001235979740865918
    "T20" = 0.00001235979740865918
    "G21" = 0.693788
    "H21" = 2.081364
    "I21" = 0.0005094032984297506
    "J21" = 0.0005094032984297505
    "M21" = 0.6405483333333334
    "N21" = 1.921645
    "O21" = 0.11614513789011
    "P21" = 0.11614513789011
    "Q21" = 0.4444047470866667
    "R21" = 3.99964272378
    "S21" = 0.00005916471633780026
    "T21" = 0.00005916471633780024
    "G22" = 0.1593103333333333
    "H22" = 0.477931
    "I22" = 0.0001169711918827409
    "J22" = 0.0001169711918827409
    "M22" = 3.339352
    "N22" = 10.018056
    "O22" = 0.6054960700393903
    "P22" = 0.6054960700393903
    "Q22" = 0.5319932802373334
    "R22" = 4.787939522136001
    "S22" = 0.00007082559699282306
    "T22" = 0.00007082559699282304
    "G23" = 0.1593103333333333
    "H23" = 0.477931
    "I23" = 0.0001169711918827409
    "J23" = 0.0001169711918827409
    "O23" = 0.2540955070726236
    "P23" = 0.2540955070726236
    "Q23" = 0.223250172858
    "R23" = 2.009251555722
    "S23" = 0.00002972185431433421
    "T23" = 0.0000297218543143342
    "G24" = 0.1593103333333333
    "H24" = 0.477931
    "I24" = 0.0001169711918827409
    "J24" = 0.0001169711918827409
    "K24" = 2
    "L24" = 0.6666666666666666
    "M24" = 0.1338136666666667
    "N24" = 0.401441
    "O24" = 0.02426328499787613
    "P24" = 0.02426328499787612
    "Q24" = 0.02131789984122223
    "R24" = 0.191861098571
    "S24" = 0.000002838105365192197
    "T24" = 0.000002838105365192196
    "G25" = 0.1593103333333333
    "H25" = 0.477931
    "I25" = 0.0001169711918827409
    "J25" = 0.0001169711918827409
    "M25" = 0.6405483333333334
    "N25" = 1.921645
    "O25" = 0.11614513789011
    "P25" = 0.11614513789011
    "Q25" = 0.1020459684994444
    "R25" = 0.918413716495
    "S25" = 0.00001358563521039146
    "T25" = 0.00001358563521039146
}

foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

Write-Host "Applied $($changes.Count) cell updates"
